# Apply the edit described by the diff to slide 7:
#  - Merge the two runs of "The XML response is not written " + "using attributes"
#    into a single run.
#  - Add a blank paragraph, a "Link to the GitHub Repository:" paragraph, a
#    hyperlinked URL paragraph, and a trailing blank paragraph.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Step 1: merge the "The XML response is not written " / "using attributes" runs ---
$idx = $tr.Text.IndexOf("The XML response is not written ")
$len = $tr.Length - $idx
$paraRange = $tr.Characters($idx + 1, $len)
$paraRange.Text = "The XML response is not written using attributes"

# --- Step 2: append the new paragraphs ---
# New paragraph separator is a carriage return; new paragraphs inherit the
# indent level (lvl="1") of the paragraph they follow, so we fix up levels
# for the paragraphs that need to differ afterwards.
$newText = "`r`rLink to the GitHub Repository:`rhttps://github.com/HyderickCSarrell/CS488_APIDOCPRESENTATION`r"
$inserted = $tr.InsertAfter($newText)

# Refresh full text range after insertion
$tr = $sh.TextFrame.TextRange

# --- Step 3: fix up indent levels ---
# "Link to the GitHub Repository:" paragraph should be at the top indent level (IndentLevel=1 -> lvl omitted/0)
$linkLabelStart = $tr.Text.IndexOf("Link to the GitHub Repository:") + 1
$linkLabelLen = "Link to the GitHub Repository:".Length
$tr.Characters($linkLabelStart, $linkLabelLen).IndentLevel = 1

# --- Step 4: add the hyperlink to the URL run ---
$url = "https://github.com/HyderickCSarrell/CS488_APIDOCPRESENTATION"
$urlStart = $tr.Text.IndexOf($url) + 1
$urlLen = $url.Length
$urlRange = $tr.Characters($urlStart, $urlLen)
$urlRange.ActionSettings(1).Hyperlink.Address = $url
